$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "ExpectedFilenames" column (K) values for rows 2-13 to reflect the
# new Standard/Excel/Word report naming convention introduced in this revision.
$ws.Range("K2").Value = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Clinical-2023_"
$ws.Range("K3").Value = "ExcelReport-NewImportLogic_1-Test_Automation_1-Clinical-"
$ws.Range("K4").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Clinical-"
$ws.Range("K5").Value = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Economic-2023_"
$ws.Range("K6").Value = "ExcelReport-NewImportLogic_1-Test_Automation_1-Economic-"
$ws.Range("K7").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Economic-"
$ws.Range("K8").Value = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Quality of Life-2023_"
$ws.Range("K9").Value = "ExcelReport-NewImportLogic_1-Test_Automation_1-Quality of Life-"
$ws.Range("K10").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Quality of Life-"
$ws.Range("K11").Value = "StandardExcelReport-NewImportLogic_1 - Test_Automation_1-Real-world Evidence-2023_"
$ws.Range("K12").Value = "ExcelReport-NewImportLogic_1-Test_Automation_1-Real-world Evidence-"
$ws.Range("K13").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Real-world Evidence-"

# Update the sheet view: scroll so column H is first visible, and select K2:K13
# with K2 as the active cell.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("K2:K13").Select()
